$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Actual Production (MW)" values for rows 2..97 (96 quarter-hour slots),
# refreshed to reflect the newly fetched day's data.
$newB = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,5,19,40,71,103,140,183,253,316,345,387,462,543,606,661,672,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($r = 2; $r -le 97; $r++) {
    $idx = $r - 2

    # Shift the timestamp forward by 3 days (the data pull now covers a later day).
    $curDate = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $curDate + 3

    # Replace the production figure with the refreshed value.
    $ws.Cells.Item($r, 2).Value = $newB[$idx]
}
